# Apply "Week 15 simulations" update to Target Depth Data.xlsx
# Row 3 (label "R") on both the OFF and DEF sheets gets updated values
# for columns B (Short Att) through F (Short Int). G (Deep Int) is unchanged.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 510
$wsOff.Range("C3").Value = 363
$wsOff.Range("D3").Value = 123
$wsOff.Range("E3").Value = 71
$wsOff.Range("F3").Value = 10

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 507
$wsDef.Range("C3").Value = 329
$wsDef.Range("D3").Value = 108
$wsDef.Range("E3").Value = 53
